$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "View Module (Licanthrope)" -> "View Module (Lican" + ")"  (two runs)
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("View Module (Licanthrope)", $false, $false, $false, $false, $false, $true, 1, $false, "View Module (Lican)", 2)

# Force the trailing ")" onto its own run: briefly drop a collapsed
# bookmark right before it (adding a bookmark splits the surrounding run at
# that point) and then delete the bookmark again - the run split survives.
$r = $d.Content
$r.Find.Execute("View Module (Lican", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint = $d.Range($r.End, $r.End)
$d.Bookmarks.Add("TmpSplit1", $splitPoint) | Out-Null
$d.Bookmarks("TmpSplit1").Delete()

# ---------------------------------------------------------------------------
# 2) "Database module()" -> "Database module(" + "BloodLine" + ")" (three runs)
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("Database module()", $true, $false, $false, $false, $false, $true, 1, $false, "Database module(BloodLine)", 2)

$r = $d.Content
$r.Find.Execute("BloodLine", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$beforeBloodLine = $d.Range($r.Start, $r.Start)
$d.Bookmarks.Add("TmpSplit2a", $beforeBloodLine) | Out-Null
$d.Bookmarks("TmpSplit2a").Delete()
$afterBloodLine = $d.Range($r.End, $r.End)
$d.Bookmarks.Add("TmpSplit2b", $afterBloodLine) | Out-Null
$d.Bookmarks("TmpSplit2b").Delete()

# ---------------------------------------------------------------------------
# 3) "Batch module()" -> "Batch module(" + "Odin" + [[_GoBack bookmark]] + ")"
#    The _GoBack bookmark is relocated here from near the end of the document.
# ---------------------------------------------------------------------------

# Remove the old (now-empty) _GoBack bookmark first so the name is free.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$r = $d.Content
$r.Find.Execute("Batch module()", $true, $false, $false, $false, $false, $true, 1, $false, "Batch module(Odin)", 2)

$r = $d.Content
$r.Find.Execute("Odin", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$beforeOdin = $d.Range($r.Start, $r.Start)
$d.Bookmarks.Add("TmpSplit3", $beforeOdin) | Out-Null
$d.Bookmarks("TmpSplit3").Delete()

$r = $d.Content
$r.Find.Execute("Odin", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterOdin = $d.Range($r.End, $r.End)
$d.Bookmarks.Add("_GoBack", $afterOdin) | Out-Null
